$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2070552147239264
$ws.Range("C2").Value = 0.5383435582822086
$ws.Range("J2").Value = 0.01073619631901841
$ws.Range("P2").Value = 0.1595092024539877
$ws.Range("S2").Value = 0.08435582822085889
$ws.Range("B3").Value = 0.0108991825613079
$ws.Range("C3").Value = 0.02997275204359673
$ws.Range("J3").Value = 0.02452316076294278
$ws.Range("P3").Value = 0.7438692098092643
$ws.Range("S3").Value = 0.1907356948228883
$ws.Range("J4").Value = 0.0410958904109589
$ws.Range("O4").Value = 0.0136986301369863
$ws.Range("P4").Value = 0.7397260273972602
$ws.Range("S4").Value = 0.2054794520547945
$ws.Range("B6").Value = 0.08071748878923767
$ws.Range("D6").Value = 0.008968609865470852
$ws.Range("F6").Value = 0.06278026905829596
$ws.Range("J6").Value = 0.2533632286995516
$ws.Range("O6").Value = 0.01345291479820628
$ws.Range("Q6").Value = 0.1704035874439462
$ws.Range("R6").Value = 0.05829596412556054
$ws.Range("S6").Value = 0.352017937219731
$ws.Range("B7").Value = 0.1092896174863388
$ws.Range("D7").Value = 0.02185792349726776
$ws.Range("F7").Value = 0.04098360655737705
$ws.Range("J7").Value = 0.1229508196721311
$ws.Range("O7").Value = 0.01912568306010929
$ws.Range("Q7").Value = 0.1693989071038251
$ws.Range("R7").Value = 0.06010928961748634
$ws.Range("S7").Value = 0.4562841530054645
$ws.Range("B8").Value = 0.1138790035587189
$ws.Range("D8").Value = 0.01897983392645314
$ws.Range("F8").Value = 0.05931198102016608
$ws.Range("J8").Value = 0.129300118623962
$ws.Range("O8").Value = 0.02016607354685647
$ws.Range("Q8").Value = 0.1553973902728351
$ws.Range("R8").Value = 0.08659549228944247
$ws.Range("S8").Value = 0.4163701067615658
$ws.Range("B9").Value = 0.1181318681318681
$ws.Range("D9").Value = 0.01098901098901099
$ws.Range("F9").Value = 0.04395604395604396
$ws.Range("J9").Value = 0.1043956043956044
$ws.Range("O9").Value = 0.02472527472527472
$ws.Range("Q9").Value = 0.1785714285714286
$ws.Range("R9").Value = 0.09615384615384616
$ws.Range("S9").Value = 0.4230769230769231
$ws.Range("B10").Value = 0.1240875912408759
$ws.Range("D10").Value = 0.0176041219407471
$ws.Range("E10").Value = 0.0008587376556462001
$ws.Range("F10").Value = 0.0790038643194504
$ws.Range("J10").Value = 0.1086303134392443
$ws.Range("O10").Value = 0.0279089738085015
$ws.Range("Q10").Value = 0.2091026191498497
$ws.Range("R10").Value = 0.0734220695577501
$ws.Range("S10").Value = 0.3593817088879347
$ws.Range("G11").Value = 0.1290322580645161
$ws.Range("J11").Value = 0.1129032258064516
$ws.Range("K11").Value = 0.1790322580645161
$ws.Range("L11").Value = 0.5596774193548387
$ws.Range("S11").Value = 0.01935483870967742
$ws.Range("G12").Value = 0.7219101123595506
$ws.Range("J12").Value = 0.1853932584269663
$ws.Range("K12").Value = 0.01966292134831461
$ws.Range("L12").Value = 0.01404494382022472
$ws.Range("S12").Value = 0.05898876404494382
$ws.Range("G13").Value = 0.6811594202898551
$ws.Range("J13").Value = 0.2318840579710145
$ws.Range("S13").Value = 0.08695652173913043
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01357466063348416
$ws.Range("H15").Value = 0.1900452488687783
$ws.Range("I15").Value = 0.06561085972850679
$ws.Range("J15").Value = 0.3122171945701357
$ws.Range("K15").Value = 0.05429864253393665
$ws.Range("M15").Value = 0.002262443438914027
$ws.Range("O15").Value = 0.08597285067873303
$ws.Range("S15").Value = 0.2760180995475113
$ws.Range("F16").Value = 0.02369668246445497
$ws.Range("H16").Value = 0.1658767772511848
$ws.Range("I16").Value = 0.08767772511848342
$ws.Range("J16").Value = 0.3886255924170616
$ws.Range("K16").Value = 0.1137440758293839
$ws.Range("M16").Value = 0.01895734597156398
$ws.Range("O16").Value = 0.04265402843601896
$ws.Range("S16").Value = 0.1587677725118483
$ws.Range("F17").Value = 0.02098765432098765
$ws.Range("H17").Value = 0.154320987654321
$ws.Range("I17").Value = 0.08518518518518518
$ws.Range("J17").Value = 0.408641975308642
$ws.Range("K17").Value = 0.1172839506172839
$ws.Range("M17").Value = 0.01604938271604938
$ws.Range("O17").Value = 0.05802469135802469
$ws.Range("S17").Value = 0.1395061728395062
$ws.Range("F18").Value = 0.00911854103343465
$ws.Range("H18").Value = 0.1793313069908815
$ws.Range("I18").Value = 0.09422492401215805
$ws.Range("J18").Value = 0.3768996960486322
$ws.Range("K18").Value = 0.1337386018237082
$ws.Range("M18").Value = 0.00911854103343465
$ws.Range("O18").Value = 0.05167173252279635
$ws.Range("S18").Value = 0.1458966565349544
$ws.Range("F19").Value = 0.02830188679245283
$ws.Range("H19").Value = 0.2028301886792453
$ws.Range("I19").Value = 0.07861635220125786
$ws.Range("J19").Value = 0.345125786163522
$ws.Range("K19").Value = 0.1080974842767296
$ws.Range("M19").Value = 0.01886792452830189
$ws.Range("N19").Value = 0.0003930817610062893
$ws.Range("O19").Value = 0.06525157232704402
$ws.Range("S19").Value = 0.1525157232704402
